# Weekly price-sheet update: a new week's observation is inserted as the
# new row 2 (most-recent-first ordering), pushing all existing data rows
# down by one. The sheet's used range grows from A1:R89 to A1:R90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (the first data row),
# shifting rows 2:89 down to 3:90.
$ws.Rows.Item(2).Insert(-4121, 1)

# Excel's Insert copies the formatting of the row above (the bold header
# row) into the freshly inserted row; strip that so the new data row
# matches the plain formatting used by every other data row.
$ws.Range("A2:R2").ClearFormats()

# Populate the new row with this week's observation (same market /
# product / classification as every other row in this sheet).
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44922
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 100114007
$ws.Cells.Item(2, 7).Value = "Jengibre"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 14000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 14500
$ws.Cells.Item(2, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(2, 15).Value = "Perú"
$ws.Cells.Item(2, 16).Value = 1115
$ws.Cells.Item(2, 17).Value = 13
$ws.Cells.Item(2, 18).Value = "Hortaliza"
